# "add new pack for WNRS" - adds a new "Existential Crisis Edition" sheet
# at the end of the workbook, populated like the other "Edition" decks.

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up the view state of two existing sheets (best effort; the
#        headless runtime does not expose every window-chrome knob, but we
#        still issue the idiomatic COM calls for the parts that are wired
#        through to the exporter). --------------------------------------

# "Voting Edition" (sheet 13): scroll the viewport down so row 10 is the
# first visible row.
$votingSheet = $wb.Worksheets.Item("Voting Edition")
$votingSheet.Activate()
$excel.ActiveWindow.ScrollRow = 10

# "Self-Love Edition" (sheet 5): the cursor moves off of this sheet (it was
# previously the active tab) and A1:C5 becomes the lingering selection.
$selfLoveSheet = $wb.Worksheets.Item("Self-Love Edition")
$selfLoveSheet.Activate()
$selfLoveSheet.Range("A1:C5").Select()

# --- 2. Create the new worksheet at the very end of the workbook. -------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Existential Crisis Edition"

# --- 3. Header block (Type / Instruction / Back) -------------------------

$newSheet.Range("A1").Value = "Type"
$newSheet.Range("B1").Value = "Single-Player"

$newSheet.Range("A2").Value = "Instruction"
$newSheet.Range("B2").Value = "These are just some thought-provoking questions I encountered online and also`nsome daily reminder that might be useful as a timely reminder for yourself.`nReady?"
$newSheet.Range("B2").Orientation = 0

$newSheet.Range("A3").Value = "Back"
$newSheet.Range("B3").Value = "Might cause you to question everything you know"

# --- 4. Card / Prompt table ----------------------------------------------

$newSheet.Range("A5").Value = "Card"
$newSheet.Range("B5").Value = "Prompt"

$newSheet.Range("A6").Value = "Y"
$newSheet.Range("B6").Value = "Is there someone you’re currently blaming or angry with that you could forgive if tomorrow was your last day on earth?"

$newSheet.Range("A7").Value = "Y"
$newSheet.Range("B7").Value = "What are you most passionate about?"

$newSheet.Range("A8").Value = "Y"
$newSheet.Range("B8").Value = "Are you holding back your skills, talents, ideas, or expression so that others won’t feel intimidated?"

$newSheet.Range("A9").Value = "Y"
$newSheet.Range("B9").Value = "Are you creating your life in a way that feels joyous and inspiring?"

$newSheet.Range("A10").Value = "Y"
$newSheet.Range("B10").Value = "Do you feel powerful and confident, regardless of your accomplishments or what others think of you?"

$newSheet.Range("A11").Value = "Y"
$newSheet.Range("B11").Value = "Do you truly love who you authentically are?"

$newSheet.Range("A12").Value = "Y"
$newSheet.Range("B12").Value = "Am I feeling passionate about my goals and is what I’m doing stimulating me?"

$newSheet.Range("A13").Value = "Y"
$newSheet.Range("B13").Value = "Am I making choices from a sense of empowerment?"

$newSheet.Range("A14").Value = "Y"
$newSheet.Range("B14").Value = "Can I let go of something that isn’t working?"

$newSheet.Range("A15").Value = "Y"
$newSheet.Range("B15").Value = "Can I allow myself to feel fully whatever is coming up?"

$newSheet.Range("A16").Value = "N"
$newSheet.Range("B16").Value = "Reminder Remember you have full power over yourself, your reactions, and where you direct your focus."
